# Modifications due to how-to 549.
#
# - Switch the active/selected sheet from "Producers" to "Contracts".
# - On the "Contracts" sheet:
#     * Column C is no longer best-fit/auto-sized; it gets a fixed custom
#       width (~15.86 characters) and the selected cell becomes C11.
#     * Several "Minimal Number of Contributors" values change.

$wb = $excel.ActiveWorkbook
$producers = $wb.Worksheets.Item("Producers")
$contracts = $wb.Worksheets.Item("Contracts")

# Update the contributor counts in column C.
$contracts.Range("C2").Value = 4
$contracts.Range("C3").Value = 0
$contracts.Range("C5").Value = 0
$contracts.Range("C6").Value = 2
$contracts.Range("C7").Value = 0
$contracts.Range("C8").Value = 0
$contracts.Range("C10").Value = 2

# Column C switches from auto (best-fit) width to an explicit custom width.
$contracts.Columns.Item(3).ColumnWidth = 15

# Make "Contracts" the active sheet/tab (this clears tabSelected on
# "Producers" and sets it on "Contracts"), then leave the selection on C11
# as recorded in the saved view.
$contracts.Activate()
$contracts.Range("C11").Select()
